# Append 45 new data rows (102-146) to the master-reg_center_machine_device_h
# worksheet, continuing the existing repeating pattern of reference/lookup
# values (A cycles 10002-10010, B cycles 10021-10029, C increments by 1),
# then leave the selection on the row below the new data and switch the
# sheet's print orientation to portrait, mirroring the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuples of (A, B, C) for rows 102..146 - D/E/F/G/H are constant for every
# row and match the values already used by rows 2-101.
$data = @(
    ,(10002, 10021, 3000121)
    ,(10003, 10022, 3000122)
    ,(10004, 10023, 3000123)
    ,(10005, 10024, 3000124)
    ,(10006, 10025, 3000125)
    ,(10007, 10026, 3000126)
    ,(10008, 10027, 3000127)
    ,(10009, 10028, 3000128)
    ,(10010, 10029, 3000129)
    ,(10002, 10021, 3000130)
    ,(10003, 10022, 3000131)
    ,(10004, 10023, 3000132)
    ,(10005, 10024, 3000133)
    ,(10006, 10025, 3000134)
    ,(10007, 10026, 3000135)
    ,(10008, 10027, 3000136)
    ,(10009, 10028, 3000137)
    ,(10010, 10029, 3000138)
    ,(10002, 10021, 3000139)
    ,(10003, 10022, 3000140)
    ,(10004, 10023, 3000141)
    ,(10005, 10024, 3000142)
    ,(10006, 10025, 3000143)
    ,(10007, 10026, 3000144)
    ,(10008, 10027, 3000145)
    ,(10009, 10028, 3000146)
    ,(10010, 10029, 3000147)
    ,(10002, 10021, 3000148)
    ,(10003, 10022, 3000149)
    ,(10004, 10023, 3000150)
    ,(10005, 10024, 3000151)
    ,(10006, 10025, 3000152)
    ,(10007, 10026, 3000153)
    ,(10008, 10027, 3000154)
    ,(10009, 10028, 3000155)
    ,(10010, 10029, 3000156)
    ,(10002, 10021, 3000157)
    ,(10003, 10022, 3000158)
    ,(10004, 10023, 3000159)
    ,(10005, 10024, 3000160)
    ,(10006, 10025, 3000161)
    ,(10007, 10026, 3000162)
    ,(10008, 10027, 3000163)
    ,(10009, 10028, 3000164)
    ,(10010, 10029, 3000165)
)

$startRow = 102
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin()"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Mirror the saved selection: the cell just below the new data, through the
# bottom of the sheet.
$lastRow = $startRow + $data.Count
$ws.Range("A" + $lastRow + ":XFD1048576").Select()

# Switch the page to portrait orientation (adds <pageSetup orientation="portrait".../>)
$ws.PageSetup.Orientation = 1
